$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "11_02_2024"
$ws.Range("H2").Value = 935
$ws.Range("H3").Value = 828
$ws.Range("H4").Value = 1388
$ws.Range("H5").Value = 2869

$ws.Range("H5").Select()
